$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (shifting old row 14 data down to row 15),
# then fill new row 14 with the updated weekly data.
$ws.Rows.Item(14).Insert()

# Set row 14 (new data row) values
$ws.Cells.Item(14, 1).Value = 7
$ws.Cells.Item(14, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(14, 3).Value = "Ñuble"
$ws.Cells.Item(14, 4).Value = 44468
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(13, 4).NumberFormat
$ws.Cells.Item(14, 5).Value = 16
$ws.Cells.Item(14, 6).Value = 100112013
$ws.Cells.Item(14, 7).Value = "Alcachofa"
$ws.Cells.Item(14, 8).Value = "Madrigal"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 60
$ws.Cells.Item(14, 11).Value = 12000
$ws.Cells.Item(14, 12).Value = 13000
$ws.Cells.Item(14, 13).Value = 12500
$ws.Cells.Item(14, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(14, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(14, 16).Value = 312
$ws.Cells.Item(14, 17).Value = 40
$ws.Cells.Item(14, 18).Value = "Hortaliza"
